$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.51144003868103
$ws.Range("B1").Value = 3.165001392364502
$ws.Range("C1").Value = 2.664785385131836
$ws.Range("D1").Value = 2.055957555770874
$ws.Range("E1").Value = 1.253536820411682
